# Applies the LOB1218 "Ativação 2022" content revision:
#  - Ativação date bumped from 01/01/2020 to 01/01/2022
#  - Programa resumido / Short syllabus: the "river in natural regime" sentence removed
#  - Programa / Syllabus: rewritten as bullet-style topic lists
#  - Método / Critério: evaluation text reworded (provas -> trabalhos/exercícios)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Ativação: (row 8) -------------------------------------------------
# Leading apostrophe forces text storage so "01/01/2022" stays a literal
# string (matching the original cell's type) instead of being auto-parsed
# into a date serial number.
$ws.Range("B8").Value = "'01/01/2022"
$ws.Range("C8").Value = "'01/01/2022"

# --- Programa resumido: (row 14) ---------------------------------------
$ws.Range("B14").Value = "Barragens e Reservatórios. Usos da água demandados para o interesse humano e Panorama Geral da Engenharia dos Recursos Hídricos. Impactos Ambientais dos Usos da Água. Gestão dos Recursos Hídricos."
$ws.Range("C14").Value = "Barragens e Reservatórios. Usos da água demandados para o interesse humano e Panorama Geral da Engenharia dos Recursos Hídricos. Impactos Ambientais dos Usos da Água. Gestão dos Recursos Hídricos."

# --- Short syllabus: (row 15) -------------------------------------------
$ws.Range("B15").Value = "Barrages and Water Tanks. The use of the water demanded for the human interest and Hydric Resources Engineering General View. The Usage of the Water Enviropnmetal Impacts. Hydric Resources Management."
$ws.Range("C15").Value = "Barrages and Water Tanks. The use of the water demanded for the human interest and Hydric Resources Engineering General View. The Usage of the Water Enviropnmetal Impacts. Hydric Resources Management."

# --- Programa: (row 16) ---------------------------------------------
$ws.Range("B16").Value = "- Políticas Públicas, Balanço Hídrico,- Demanda de água e disponibilidade dos recursos hídricos: Abastecimento Humano, águas para Agropecuária e indústria. - Hidreletricidade. - Barragens e Reservatórios,- Navegação Interior.- Águas Subterrâneas.- Gerenciamento dos Recursos Hídricos.- Hidroeconomia- Relação entre saneamento e qualidade da água"
$ws.Range("C16").Value = "- Políticas Públicas, Balanço Hídrico,- Demanda de água e disponibilidade dos recursos hídricos: Abastecimento Humano, águas para Agropecuária e indústria. - Hidreletricidade. - Barragens e Reservatórios,- Navegação Interior.- Águas Subterrâneas.- Gerenciamento dos Recursos Hídricos.- Hidroeconomia- Relação entre saneamento e qualidade da água"

# --- Syllabus: (row 17) -----------------------------------------------
$ws.Range("B17").Value = "- Public Policies, Water Balance,- Water demand and availability of water resources: Human Supply, water for Agriculture and industry.- Hydroelectricity.- Dams and Reservoirs,- Inland navigation.- Groundwater.- Water Resources Management.- Hydroeconomics- Relationship between sanitation and water quality"
$ws.Range("C17").Value = "- Public Policies, Water Balance,- Water demand and availability of water resources: Human Supply, water for Agriculture and industry.- Hydroelectricity.- Dams and Reservoirs,- Inland navigation.- Groundwater.- Water Resources Management.- Hydroeconomics- Relationship between sanitation and water quality"

# --- Método: (row 19) -----------------------------------------------
$ws.Range("B19").Value = "Avaliação baseada em trabalhos com dados reais, exercícios, trabalhos práticos e relatórios."
$ws.Range("C19").Value = "Avaliação baseada em trabalhos com dados reais, exercícios, trabalhos práticos e relatórios."

# --- Critério: (row 20) -----------------------------------------------
$ws.Range("B20").Value = "Média ponderada das notas atribuídas aos exercícios e trabalhos práticos e relatórios."
$ws.Range("C20").Value = "Média ponderada das notas atribuídas aos exercícios e trabalhos práticos e relatórios."
